{"js": "// Apply the \"Built site for gh-pages\" schedule-table edits.\n// The first table on the page lists weeks in rows (row 0 = header row).\n// We update four cells (Topic/Tools columns) across weeks 2-5, replacing\n// each cell's full text while preserving the existing paragraph\n// formatting (pStyle \"Compact\") and run.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Each entry: [rowIndex, colIndex, oldText, newText]\nconst edits = [\n  [2, 2, \"Filesystems; Markup Languages; Quarto\", \"Markup Languages; Quarto\"],\n  [3, 1, \"A Field Guide to Data\", \"File Management & Version Control\"],\n  [3, 2, \"Version Control; Git; Data Formats\", \"Filesystems; Git; GitHub\"],\n  [4, 1, \"Wrangling Tidy Data\", \"A Field Guide to Data\"],\n  [4, 2, \"Tidyverse\", \"Data Formats; Tidyverse\"],\n  [5, 1, \"Flow Control\", \"Wrangling Tidy Data; Flow Control\"],\n];\n\nfor (const [rowIndex, colIndex, oldText, newText] of edits) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const results = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\n      `Could not find expected text \"${oldText}\" in row ${rowIndex}, column ${colIndex}.`\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply the \"Built site for gh-pages\" schedule-table edits.\n# The first table on the page lists weeks in rows (row 1 = header row in\n# Word's 1-based COM indexing). We update four cells (Topic/Tools columns)\n# across weeks 2-5, replacing each cell's text while preserving the\n# existing paragraph formatting (pStyle \"Compact\") and run, since setting\n# Range.Text only replaces the text up to the cell-end mark.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Each entry: row, column, expected old text, new text (1-based row/col)\n$edits = @(\n    @{ Row = 3; Col = 3; Old = \"Filesystems; Markup Languages; Quarto\"; New = \"Markup Languages; Quarto\" },\n    @{ Row = 4; Col = 2; Old = \"A Field Guide to Data\"; New = \"File Management & Version Control\" },\n    @{ Row = 4; Col = 3; Old = \"Version Control; Git; Data Formats\"; New = \"Filesystems; Git; GitHub\" },\n    @{ Row = 5; Col = 2; Old = \"Wrangling Tidy Data\"; New = \"A Field Guide to Data\" },\n    @{ Row = 5; Col = 3; Old = \"Tidyverse\"; New = \"Data Formats; Tidyverse\" },\n    @{ Row = 6; Col = 2; Old = \"Flow Control\"; New = \"Wrangling Tidy Data; Flow Control\" }\n)\n\nforeach ($edit in $edits) {\n    $cell = $t.Cell($edit.Row, $edit.Col)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($current -ne $edit.Old) {\n        throw \"Cell ($($edit.Row), $($edit.Col)) text was '$current', expected '$($edit.Old)'.\"\n    }\n\n    $cell.Range.Text = $edit.New\n}\n"}
